$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update refresh: row 3 now carries MetLife's data, row 4 now carries AIG's
# data (the underlying feed re-sorted the two issuers), and every row's
# numeric columns (price/RSI/5d return/scores/final-score) were refreshed.

# Row 2 - UnitedHealth (UNH)
$ws.Range("D2").Value = 329.26
$ws.Range("E2").Value = 55.4
$ws.Range("F2").Value = -0.15
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 46
$ws.Range("K2").Value = 54.9
$ws.Range("N2").Value = 49.16024380385575

# Row 3 - now MetLife (MET)
$ws.Range("B3").Value = "MetLife, Inc."
$ws.Range("C3").Value = "MET"
$ws.Range("D3").Value = 78.34999999999999
$ws.Range("E3").Value = 48.5
$ws.Range("F3").Value = 2.34
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 33
$ws.Range("K3").Value = 52.7
$ws.Range("N3").Value = 49.16024380385575

# Row 4 - now American International Group (AIG)
$ws.Range("B4").Value = "American International Group, I"
$ws.Range("C4").Value = "AIG"
$ws.Range("D4").Value = 76.59
$ws.Range("E4").Value = 42.2
$ws.Range("F4").Value = 0.5600000000000001
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 46
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 49.7
$ws.Range("N4").Value = 49.16024380385575

# Row 5 - Prudential Financial (PRU)
$ws.Range("D5").Value = 111.3
$ws.Range("E5").Value = 69.40000000000001
$ws.Range("F5").Value = 2.82
$ws.Range("K5").Value = 47.1
$ws.Range("N5").Value = 49.16024380385575
